$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New marriage-index rows (Linyola, 1761-1765) to append after the existing
# data (previously ending at row 767).  Each tuple is:
# (row, Any, Fotograma, CognomsFamilia, Projecte, Rollo, Serie, Anys)
$rows = @(
    @(768, 1761, 31, 'Roige Pons', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(769, 1761, 31, 'Mas Planes', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(770, 1761, 32, 'Repetit', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(771, 1761, 33, 'Valles Massot', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(772, 1761, 33, 'Segura Palou', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(773, 1761, 34, 'Druet Pujades', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(774, 1762, 34, 'Domingo Pera', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(775, 1762, 35, 'Solsona Melé', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(776, 1762, 35, 'Roma Redon', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(777, 1762, 35, 'Vallés Fabregat', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(778, 1763, 36, 'Pujol Pedrós', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(779, 1763, 36, 'Batalla Margall', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(780, 1763, 37, 'Cascalló Pujades', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(781, 1763, 37, 'Manyach Mas', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(782, 1763, 38, 'Torrà Eroles', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(783, 1763, 38, 'Vergé Codina', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(784, 1763, 39, 'Viladebaix Solà', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(785, 1763, 39, 'Vergé Codina', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(786, 1764, 40, 'Claverol Siriols', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(787, 1764, 40, 'Tarragó Solà', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(788, 1764, 40, 'Palou Vergé', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(789, 1764, 41, 'Sabater Valles', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(790, 1764, 41, 'Mas Thomas', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(791, 1764, 42, 'Roma Oriola', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(792, 1764, 42, 'Mata Mas', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(793, 1764, 43, 'Duart Alexandre', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(794, 1764, 43, 'Homs Cascalló', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(795, 1764, 44, 'Solsona Mas', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(796, 1765, 44, 'Navau Santesmases', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(797, 1765, 44, 'Galceran Rodon', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(798, 1765, 45, 'Vergé Thomas', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(799, 1765, 45, 'Roma Vallés', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(800, 1765, 46, 'Bonjorn Domenech', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(801, 1765, 46, 'Gine Serra', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(802, 1765, 46, 'Claverol Roca', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(803, 1765, 47, 'Carrera Galceran', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(804, 1765, 47, 'Cases Mas', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(805, 1765, 48, 'Gispert Agulló', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(806, 1761, 48, 'Pilfort Teixidó', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(807, 1761, 48, 'Sunyé Fabregat', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(808, 1761, 49, 'Mas Giné', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(809, 1761, 49, 'Pujades Porta', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
    @(810, 1761, 49, 'Codol ?', 'SPN 2,02 C', 47, 'A,2', '1749-1770'),
)

foreach ($row in $rows) {
    $r   = $row[0]
    $ws.Cells.Item($r, 1).Value  = $row[1]   # A - Any
    $ws.Cells.Item($r, 2).Value  = $row[2]   # B - Fotograma
    $ws.Cells.Item($r, 3).Value  = $row[3]   # C - Cognoms Familia
    $ws.Cells.Item($r, 10).Value = $row[4]   # J - Projecte
    $ws.Cells.Item($r, 11).Value = $row[5]   # K - Rollo
    $ws.Cells.Item($r, 12).Value = $row[6]   # L - Serie
    $ws.Cells.Item($r, 13).Value = $row[7]   # M - Anys
}

# Move the active selection to the next empty row, as Excel does after
# manual data entry.
$ws.Range("A811").Select()
